$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$wsALC.Range("H17").Value = 4251.2856
$wsALC.Range("I17").Value = 8563
$wsALC.Range("J17").Value = 3532.6667
$wsALC.Range("K17").Value = 25689
$wsALC.Range("L17").Value = 10598.0001
$wsALC.Range("M17").Value = -25521
$wsALC.Range("N17").Value = -10934.0001

# Row 121 (ALC)
$wsALC.Range("H121").Value = 868.36365
$wsALC.Range("I121").Value = 516
$wsALC.Range("J121").Value = 1000.5
$wsALC.Range("K121").Value = 1548
$wsALC.Range("L121").Value = 3001.5
$wsALC.Range("M121").Value = 199
$wsALC.Range("N121").Value = -6495.5

# Row 135 (ALC)
$wsALC.Range("H135").Value = 828.375
$wsALC.Range("I135").Value = 622.2
$wsALC.Range("J135").Value = 1172
$wsALC.Range("K135").Value = 5599.8
$wsALC.Range("L135").Value = 10548
$wsALC.Range("M135").Value = -3064.8
$wsALC.Range("N135").Value = -15618

# Row 137 (ALC)
$wsALC.Range("H137").Value = 60305.293
$wsALC.Range("I137").Value = 1475
$wsALC.Range("J137").Value = 78406.92
$wsALC.Range("K137").Value = 4425
$wsALC.Range("L137").Value = 235220.76
$wsALC.Range("M137").Value = -1875
$wsALC.Range("N137").Value = -240320.76

# Row 141 (ALC)
$wsALC.Range("H141").Value = 4669624.5
$wsALC.Range("I141").Value = 9334064
$wsALC.Range("J141").Value = 5185
$wsALC.Range("K141").Value = 28002192
$wsALC.Range("L141").Value = 15555
$wsALC.Range("M141").Value = -27997012
$wsALC.Range("N141").Value = -25915

$wsARM = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$wsARM.Range("H2").Value = 397471.22
$wsARM.Range("I2").Value = 427983.38
$wsARM.Range("J2").Value = 813
$wsARM.Range("K2").Value = 427983.38
$wsARM.Range("L2").Value = 813
$wsARM.Range("M2").Value = -427870.38
$wsARM.Range("N2").Value = -1039

# Row 15 (ARM)
$wsARM.Range("H15").Value = 0
$wsARM.Range("I15").Value = 0
$wsARM.Range("J15").Value = 0
$wsARM.Range("K15").Value = 0
$wsARM.Range("L15").Value = 0
$wsARM.Range("N15").ClearContents()

# Row 61 (ARM)
$wsARM.Range("H61").Value = 29368.9
$wsARM.Range("I61").Value = 35848.652
$wsARM.Range("J61").Value = 8078.2856
$wsARM.Range("K61").Value = 35848.652
$wsARM.Range("L61").Value = 8078.2856
$wsARM.Range("M61").Value = -35636.652
$wsARM.Range("N61").Value = -8502.285599999999

# Row 116 (ARM)
$wsARM.Range("H116").Value = 397471.22
$wsARM.Range("I116").Value = 427983.38
$wsARM.Range("J116").Value = 813
$wsARM.Range("K116").Value = 427983.38
$wsARM.Range("L116").Value = 813
$wsARM.Range("M116").Value = -425689.38
$wsARM.Range("N116").Value = -5401

# Row 119 (ARM)
$wsARM.Range("H119").Value = 31000
$wsARM.Range("I119").Value = 0
$wsARM.Range("J119").Value = 31000
$wsARM.Range("K119").Value = 0
$wsARM.Range("L119").Value = 31000
$wsARM.Range("N119").Value = -40676

# Row 123 (ARM)
$wsARM.Range("H123").Value = 82000
$wsARM.Range("I123").Value = 0
$wsARM.Range("J123").Value = 82000
$wsARM.Range("K123").Value = 0
$wsARM.Range("L123").Value = 82000
$wsARM.Range("N123").Value = -91800

# Row 136 (ARM)
$wsARM.Range("H136").Value = 29368.9
$wsARM.Range("I136").Value = 35848.652
$wsARM.Range("J136").Value = 8078.2856
$wsARM.Range("K136").Value = 107545.956
$wsARM.Range("L136").Value = 24234.8568
$wsARM.Range("M136").Value = -104995.956
$wsARM.Range("N136").Value = -29334.8568

$wsBSM = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$wsBSM.Range("H3").Value = 397471.22
$wsBSM.Range("I3").Value = 427983.38
$wsBSM.Range("J3").Value = 813
$wsBSM.Range("K3").Value = 427983.38
$wsBSM.Range("L3").Value = 813
$wsBSM.Range("M3").Value = -427869.38
$wsBSM.Range("N3").Value = -1041

# Row 86 (BSM)
$wsBSM.Range("H86").Value = 2519.3333
$wsBSM.Range("I86").Value = 1306

# Row 89 (BSM)
$wsBSM.Range("H89").Value = 2519.3333
$wsBSM.Range("I89").Value = 1306

# Row 132 (BSM)
$wsBSM.Range("H132").Value = 40647.75
$wsBSM.Range("I132").Value = 0
$wsBSM.Range("J132").Value = 40647.75
$wsBSM.Range("K132").Value = 0
$wsBSM.Range("L132").Value = 40647.75
$wsBSM.Range("N132").Value = -50767.75

# Row 134 (BSM)
$wsBSM.Range("H134").Value = 6012.9375
$wsBSM.Range("I134").Value = 6971.7393
$wsBSM.Range("J134").Value = 3562.6667
$wsBSM.Range("K134").Value = 20915.2179
$wsBSM.Range("L134").Value = 10688.0001
$wsBSM.Range("M134").Value = -18380.2179
$wsBSM.Range("N134").Value = -15758.0001

$wsCRP = $wb.Worksheets.Item("CRP")
# Row 36 (CRP)
$wsCRP.Range("H36").Value = 0
$wsCRP.Range("I36").Value = 0
$wsCRP.Range("J36").Value = 0
$wsCRP.Range("K36").Value = 0
$wsCRP.Range("L36").Value = 0
$wsCRP.Range("N36").ClearContents()

# Row 40 (CRP)
$wsCRP.Range("H40").Value = 0
$wsCRP.Range("I40").Value = 0
$wsCRP.Range("J40").Value = 0
$wsCRP.Range("K40").Value = 0
$wsCRP.Range("L40").Value = 0
$wsCRP.Range("N40").ClearContents()

# Row 86 (CRP)
$wsCRP.Range("H86").Value = 2317
$wsCRP.Range("I86").Value = 2060.125
$wsCRP.Range("J86").Value = 2610.5715
$wsCRP.Range("K86").Value = 2060.125
$wsCRP.Range("L86").Value = 2610.5715
$wsCRP.Range("M86").Value = -937.125
$wsCRP.Range("N86").Value = -4856.5715

# Row 89 (CRP)
$wsCRP.Range("H89").Value = 2317
$wsCRP.Range("I89").Value = 2060.125
$wsCRP.Range("J89").Value = 2610.5715
$wsCRP.Range("K89").Value = 10300.625
$wsCRP.Range("L89").Value = 13052.8575
$wsCRP.Range("M89").Value = -4684.625
$wsCRP.Range("N89").Value = -24284.8575

# Row 99 (CRP)
$wsCRP.Range("H99").Value = 2633.5881
$wsCRP.Range("I99").Value = 2431.2
$wsCRP.Range("J99").Value = 2922.7144
$wsCRP.Range("K99").Value = 2431.2
$wsCRP.Range("L99").Value = 2922.7144
$wsCRP.Range("M99").Value = -933.1999999999998
$wsCRP.Range("N99").Value = -5918.7144

# Row 121 (CRP)
$wsCRP.Range("H121").Value = 0
$wsCRP.Range("I121").Value = 0
$wsCRP.Range("J121").Value = 0
$wsCRP.Range("K121").Value = 0
$wsCRP.Range("L121").Value = 0
$wsCRP.Range("N121").ClearContents()

# Row 126 (CRP)
$wsCRP.Range("H126").Value = 2633.5881
$wsCRP.Range("I126").Value = 2431.2
$wsCRP.Range("J126").Value = 2922.7144
$wsCRP.Range("K126").Value = 7293.599999999999
$wsCRP.Range("L126").Value = 8768.143199999999
$wsCRP.Range("M126").Value = -4823.599999999999
$wsCRP.Range("N126").Value = -13708.1432

# Row 132 (CRP)
$wsCRP.Range("H132").Value = 1613.3043
$wsCRP.Range("I132").Value = 1181.75
$wsCRP.Range("J132").Value = 2599.7144
$wsCRP.Range("K132").Value = 3545.25
$wsCRP.Range("L132").Value = 7799.1432
$wsCRP.Range("M132").Value = -1015.25
$wsCRP.Range("N132").Value = -12859.1432

# Row 134 (CRP)
$wsCRP.Range("H134").Value = 2375.5715
$wsCRP.Range("I134").Value = 1981.6471
$wsCRP.Range("J134").Value = 4049.75
$wsCRP.Range("K134").Value = 5944.9413
$wsCRP.Range("L134").Value = 12149.25
$wsCRP.Range("M134").Value = -3409.9413
$wsCRP.Range("N134").Value = -17219.25

$wsCUL = $wb.Worksheets.Item("CUL")
# Row 11 (CUL)
$wsCUL.Range("H11").Value = 167864.5
$wsCUL.Range("I11").Value = 334066
$wsCUL.Range("J11").Value = 1663
$wsCUL.Range("K11").Value = 1002198
$wsCUL.Range("L11").Value = 4989
$wsCUL.Range("M11").Value = -1002058
$wsCUL.Range("N11").Value = -5269

# Row 50 (CUL)
$wsCUL.Range("H50").Value = 142958050
$wsCUL.Range("I50").Value = 348614.5
$wsCUL.Range("J50").Value = 200001820
$wsCUL.Range("K50").Value = 1045843.5
$wsCUL.Range("L50").Value = 600005460
$wsCUL.Range("M50").Value = -1045362.5
$wsCUL.Range("N50").Value = -600006422

# Row 53 (CUL)
$wsCUL.Range("H53").Value = 142958050
$wsCUL.Range("I53").Value = 348614.5
$wsCUL.Range("J53").Value = 200001820
$wsCUL.Range("K53").Value = 1045843.5
$wsCUL.Range("L53").Value = 600005460
$wsCUL.Range("M53").Value = -1045362.5
$wsCUL.Range("N53").Value = -600006422

# Row 57 (CUL)
$wsCUL.Range("H57").Value = 4000
$wsCUL.Range("I57").Value = 0
$wsCUL.Range("J57").Value = 4000
$wsCUL.Range("K57").Value = 0
$wsCUL.Range("L57").Value = 12000
$wsCUL.Range("N57").Value = -13118

# Row 131 (CUL)
$wsCUL.Range("H131").Value = 16301.018
$wsCUL.Range("I131").Value = 450
$wsCUL.Range("J131").Value = 16867.125
$wsCUL.Range("K131").Value = 1350
$wsCUL.Range("L131").Value = 50601.375
$wsCUL.Range("M131").Value = 3690
$wsCUL.Range("N131").Value = -60681.375

$wsGSM = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$wsGSM.Range("H80").Value = 2853.3333
$wsGSM.Range("I80").Value = 2732.8333
$wsGSM.Range("J80").Value = 3335.3333
$wsGSM.Range("K80").Value = 2732.8333
$wsGSM.Range("L80").Value = 3335.3333
$wsGSM.Range("M80").Value = -1734.8333
$wsGSM.Range("N80").Value = -5331.3333

# Row 83 (GSM)
$wsGSM.Range("H83").Value = 2853.3333
$wsGSM.Range("I83").Value = 2732.8333
$wsGSM.Range("J83").Value = 3335.3333
$wsGSM.Range("K83").Value = 13664.1665
$wsGSM.Range("L83").Value = 16676.6665
$wsGSM.Range("M83").Value = -8672.166499999999
$wsGSM.Range("N83").Value = -26660.6665

# Row 132 (GSM)
$wsGSM.Range("H132").Value = 1101652
$wsGSM.Range("I132").Value = 1376042.6
$wsGSM.Range("J132").Value = 4089.4285
$wsGSM.Range("K132").Value = 4128127.8
$wsGSM.Range("L132").Value = 12268.2855
$wsGSM.Range("M132").Value = -4125597.8
$wsGSM.Range("N132").Value = -17328.2855

$wsLTW = $wb.Worksheets.Item("LTW")
# Row 93 (LTW)
$wsLTW.Range("H93").Value = 1090.2609
$wsLTW.Range("I93").Value = 621.44446
$wsLTW.Range("J93").Value = 2778
$wsLTW.Range("K93").Value = 621.44446
$wsLTW.Range("L93").Value = 2778
$wsLTW.Range("M93").Value = 626.55554
$wsLTW.Range("N93").Value = -5274

# Row 119 (LTW)
$wsLTW.Range("H119").Value = 0
$wsLTW.Range("I119").Value = 0
$wsLTW.Range("J119").Value = 0
$wsLTW.Range("K119").Value = 0
$wsLTW.Range("L119").Value = 0
$wsLTW.Range("N119").ClearContents()

$wsWVR = $wb.Worksheets.Item("WVR")
# Row 28 (WVR)
$wsWVR.Range("H28").Value = 10000
$wsWVR.Range("I28").Value = 10000
$wsWVR.Range("J28").Value = 10000
$wsWVR.Range("K28").Value = 10000
$wsWVR.Range("L28").Value = 10000
$wsWVR.Range("M28").Value = -9652
$wsWVR.Range("N28").Value = -10696

# Row 119 (WVR)
$wsWVR.Range("H119").Value = 0
$wsWVR.Range("I119").Value = 0
$wsWVR.Range("J119").Value = 0
$wsWVR.Range("K119").Value = 0
$wsWVR.Range("L119").Value = 0
$wsWVR.Range("N119").ClearContents()
